# edit.ps1
# Applies the "Fix false person detection and address trailing text" change:
#   1. [[PERSON_2]] ("Fakultni nemocnice [[PERSON_2]]") was a false-positive
#      person detection on a place/saint name -> restore literal text
#      "Svate Markety" and renumber every following PERSON_N (N>=3) down by one.
#   2. ADDRESS_RE used to swallow the newline + following label, producing
#      "Sidlo: ... [[ADDRESS_1]]: [[ICO_1]]" and "[[ADDRESS_2]][[PHONE_2]]"
#      on one line; split those back into separate lines (w:br).

$d = $word.ActiveDocument

function Replace-Literal([string]$find, [string]$replaceWith) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceWith, 2) | Out-Null
}

# --- 1. False "person" detection: "[[PERSON_2]]" was really "Svate Markety" ---
Replace-Literal "Fakultní nemocnice [[PERSON_2]], příspěvková organizace" `
                "Fakultní nemocnice Svaté Markéty, příspěvková organizace"

# --- 2. Renumber PERSON_3 .. PERSON_29 down to PERSON_2 .. PERSON_28 ---
# Ascending order is safe: at step N we search only for the still-untouched
# literal "[[PERSON_N]]" and write "[[PERSON_(N-1)]]"; later steps search for
# higher, not-yet-processed numbers, so they can never re-match text this
# step just wrote.
for ($n = 3; $n -le 29; $n++) {
    Replace-Literal "[[PERSON_$n]]" "[[PERSON_$($n - 1)]]"
}

# --- 3. Fix the address block for the provider: drop trailing ": [[ICO_1]]"
#        swallowed onto the "Sidlo" line and restore the postal code + break ---
Replace-Literal "Sídlo: Nad Kampusem 821/4, [[ADDRESS_1]]: [[ICO_1]]" `
                 "Sídlo: [[ADDRESS_1]], 102 00^l[[ICO_1]]"

# --- 4. Fix the address block for the patient: split the run-on
#        "[[ADDRESS_2]][[PHONE_2]]" back onto two lines ---
Replace-Literal "[[ADDRESS_2]][[PHONE_2]]" "[[ADDRESS_2]]^l[[PHONE_2]]"
